# Report regeneration: update localization status text and refresh the
# "Status" column widths (auto-sized to the new, shorter text) on every
# sheet that surfaces it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status cells (columns E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# --- zh-cn sheet: Status column (column C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

# --- de-de sheet: Status column (column C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
